$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/077f86f77b8e0495246306e21c82869509901565/e2e/298232d8-2b59-4f22-a90c-b752c15be540.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/516a7a10d0abdd4108571496659fb2905f5f0660/e2e/298232d8-2b59-4f22-a90c-b752c15be540.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/516a7a10d0abdd4108571496659fb2905f5f0660/e2e/298232d8-2b59-4f22-a90c-b752c15be540.md"
$displayName = "298232d8-2b59-4f22-a90c-b752c15be540.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J7").Value = "298232d8-2b59-4f22-a90c-b752c15be540.9a73bebff7de3fec107b27fff0652463d8b39e3c.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-06 17:33:58"
$wsZh.Range("P7").Value = $errorDetail
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestUrl, "", "", $displayName)

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J7").Value = "298232d8-2b59-4f22-a90c-b752c15be540.9a73bebff7de3fec107b27fff0652463d8b39e3c.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-06 17:34:18"
$wsDe.Range("P7").Value = $errorDetail
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestUrl, "", "", $displayName)

Write-Host "done"
